# "Generate Report for Handback"
# The handback-status report was regenerated, which refreshes the
# timestamp columns ("Latest HO Xliff Generate Date",
# "Correspond Handoff Datetime" and "Correspond Handback DateTime")
# for each localized-file row.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" for the first data row
$wsOverview.Range("G2").Value = "2016-08-29 09:09:26"

# zh-cn sheet - Correspond Handoff / Handback datetimes for the first data row
$wsZhCn.Range("H2").Value = "2016-08-29 09:09:22"
$wsZhCn.Range("K2").Value = "2016-08-29 09:09:42"

# de-de sheet - Correspond Handoff / Handback datetimes for the first data row
$wsDeDe.Range("H2").Value = "2016-08-29 09:09:26"
$wsDeDe.Range("K2").Value = "2016-08-29 09:09:49"
